# Auto-generated edit script applying the Sagittarius_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1156.2222
$ws.Range("I132").Value = 1200.75
$ws.Range("K132").Value = 3602.25
$ws.Range("M132").Value = -1072.25

$ws.Range("H135").Value = 772.41174
$ws.Range("I135").Value = 700.4
$ws.Range("K135").Value = 6303.599999999999
$ws.Range("M135").Value = -3768.599999999999

$ws.Range("H137").Value = 2269
$ws.Range("I137").Value = 2269
$ws.Range("K137").Value = 6807
$ws.Range("M137").Value = -4257


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8992.237999999999
$ws.Range("I32").Value = 7191.9
$ws.Range("K32").Value = 7191.9
$ws.Range("M32").Value = -6904.9

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H109").Value = 90377
$ws.Range("J109").Value = 90377
$ws.Range("L109").Value = 90377
$ws.Range("N109").Value = -93151

$ws.Range("H112").Value = 46924.332
$ws.Range("J112").Value = 46924.332
$ws.Range("L112").Value = 46924.332
$ws.Range("N112").Value = -49878.332

$ws.Range("H124").Value = 19494.5
$ws.Range("J124").Value = 19494.5
$ws.Range("L124").Value = 19494.5
$ws.Range("N124").Value = -29314.5

$ws.Range("H132").Value = 1588.7693
$ws.Range("I132").Value = 1554.5834
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 4663.7502
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -2133.7502
$ws.Range("N132").Value = -11057


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 6868.6
$ws.Range("I94").Value = 8362.125
$ws.Range("J94").Value = 894.5
$ws.Range("K94").Value = 8362.125
$ws.Range("L94").Value = 894.5
$ws.Range("M94").Value = -7911.125
$ws.Range("N94").Value = -1796.5

$ws.Range("H134").Value = 719.6
$ws.Range("I134").Value = 499.5
$ws.Range("J134").Value = 1600
$ws.Range("K134").Value = 1498.5
$ws.Range("L134").Value = 4800
$ws.Range("M134").Value = 1036.5
$ws.Range("N134").Value = -9870


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1497.7727
$ws.Range("I16").Value = 1157.7222
$ws.Range("J16").Value = 3028
$ws.Range("K16").Value = 1157.7222
$ws.Range("L16").Value = 3028
$ws.Range("M16").Value = -870.7221999999999
$ws.Range("N16").Value = -3602

$ws.Range("H33").Value = 1874.75
$ws.Range("I33").Value = 1874.75
$ws.Range("K33").Value = 1874.75
$ws.Range("M33").Value = -1495.75

$ws.Range("H58").Value = 1936.8125
$ws.Range("J58").Value = 2060.7144
$ws.Range("L58").Value = 2060.7144
$ws.Range("N58").Value = -2466.7144

$ws.Range("H99").Value = 2628.1428
$ws.Range("I99").Value = 1999.3334
$ws.Range("J99").Value = 3099.75
$ws.Range("K99").Value = 1999.3334
$ws.Range("L99").Value = 3099.75
$ws.Range("M99").Value = -501.3334
$ws.Range("N99").Value = -6095.75

$ws.Range("H113").Value = 1497.7727
$ws.Range("I113").Value = 1157.7222
$ws.Range("J113").Value = 3028
$ws.Range("K113").Value = 1157.7222
$ws.Range("L113").Value = 3028
$ws.Range("M113").Value = 1012.2778
$ws.Range("N113").Value = -7368

$ws.Range("H122").Value = 1652.2667
$ws.Range("I122").Value = 1239.5555
$ws.Range("J122").Value = 2271.3333
$ws.Range("K122").Value = 3718.6665
$ws.Range("L122").Value = 6813.999899999999
$ws.Range("M122").Value = -1268.6665
$ws.Range("N122").Value = -11713.9999

$ws.Range("H126").Value = 2628.1428
$ws.Range("I126").Value = 1999.3334
$ws.Range("J126").Value = 3099.75
$ws.Range("K126").Value = 5998.0002
$ws.Range("L126").Value = 9299.25
$ws.Range("M126").Value = -3528.0002
$ws.Range("N126").Value = -14239.25

$ws.Range("H134").Value = 4059.6667
$ws.Range("I134").Value = 3942.125
$ws.Range("K134").Value = 11826.375
$ws.Range("M134").Value = -9291.375

$ws.Range("H136").Value = 1936.8125
$ws.Range("J136").Value = 2060.7144
$ws.Range("L136").Value = 6182.1432
$ws.Range("N136").Value = -11282.1432


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2518.6667
$ws.Range("J51").Value = 2410.7144
$ws.Range("L51").Value = 7232.1432
$ws.Range("N51").Value = -8152.1432

$ws.Range("H80").Value = 2495
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 2495
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("H86").Value = 7733.3335
$ws.Range("I86").Value = 200
$ws.Range("K86").Value = 600
$ws.Range("M86").Value = 586

$ws.Range("H89").Value = 7733.3335
$ws.Range("I89").Value = 200
$ws.Range("K89").Value = 1800
$ws.Range("M89").Value = 4128

$ws.Range("H131").Value = 2645.353
$ws.Range("I131").Value = 1591
$ws.Range("J131").Value = 3084.6667
$ws.Range("K131").Value = 4773
$ws.Range("L131").Value = 9254.000100000001
$ws.Range("M131").Value = 267
$ws.Range("N131").Value = -19334.0001

$ws.Range("H132").Value = 7490.25
$ws.Range("I132").Value = 3653.6667
$ws.Range("K132").Value = 32883.0003
$ws.Range("M132").Value = -30353.0003


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4350.375
$ws.Range("I70").Value = 4321.2
$ws.Range("K70").Value = 4321.2
$ws.Range("M70").Value = -4051.2

$ws.Range("H73").Value = 4350.375
$ws.Range("I73").Value = 4321.2
$ws.Range("K73").Value = 4321.2
$ws.Range("M73").Value = -3385.2

$ws.Range("H80").Value = 3292.9285
$ws.Range("J80").Value = 3600.3333
$ws.Range("L80").Value = 3600.3333
$ws.Range("N80").Value = -5596.3333

$ws.Range("H83").Value = 3292.9285
$ws.Range("J83").Value = 3600.3333
$ws.Range("L83").Value = 18001.6665
$ws.Range("N83").Value = -27985.6665

$ws.Range("H111").Value = 80000
$ws.Range("J111").Value = 80000
$ws.Range("L111").Value = 80000
$ws.Range("N111").Value = -86134

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H126").Value = 3547.5
$ws.Range("I126").Value = 3396.5
$ws.Range("J126").Value = 3849.5
$ws.Range("K126").Value = 10189.5
$ws.Range("L126").Value = 11548.5
$ws.Range("M126").Value = -7719.5
$ws.Range("N126").Value = -16488.5


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3410.16
$ws.Range("I61").Value = 3253.25
$ws.Range("J61").Value = 4037.8
$ws.Range("K61").Value = 3253.25
$ws.Range("L61").Value = 4037.8
$ws.Range("M61").Value = -3051.25
$ws.Range("N61").Value = -4441.8

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H110").Value = 49999
$ws.Range("J110").Value = 49999
$ws.Range("L110").Value = 49999
$ws.Range("N110").Value = -58179

$ws.Range("H113").Value = 3410.16
$ws.Range("I113").Value = 3253.25
$ws.Range("J113").Value = 4037.8
$ws.Range("K113").Value = 3253.25
$ws.Range("L113").Value = 4037.8
$ws.Range("M113").Value = -1083.25
$ws.Range("N113").Value = -8377.799999999999


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10799.333
$ws.Range("I41").Value = 11900
$ws.Range("J41").Value = 10579.2
$ws.Range("K41").Value = 11900
$ws.Range("L41").Value = 10579.2
$ws.Range("M41").Value = -11510
$ws.Range("N41").Value = -11359.2

$ws.Range("H103").Value = 35000
$ws.Range("J103").Value = 35000
$ws.Range("L103").Value = 35000
$ws.Range("N103").Value = -37344

$ws.Range("H107").Value = 1151.2354
$ws.Range("J107").Value = 1185.1333
$ws.Range("L107").Value = 3555.3999
$ws.Range("N107").Value = -7395.3999

$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H122").Value = 746.4375
$ws.Range("I122").Value = 753.26666
$ws.Range("K122").Value = 2259.79998
$ws.Range("M122").Value = 190.2000200000002

$ws.Range("H135").Value = 74500
$ws.Range("J135").Value = 74500
$ws.Range("L135").Value = 74500
$ws.Range("N135").Value = -84640

$ws.Range("H136").Value = 1996.6818
$ws.Range("I136").Value = 2090.3333
$ws.Range("J136").Value = 1575.25
$ws.Range("K136").Value = 6270.999899999999
$ws.Range("L136").Value = 4725.75
$ws.Range("M136").Value = -3720.999899999999
$ws.Range("N136").Value = -9825.75

